$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 37040788
$ws.Range("I98").Value = 45458476
$ws.Range("K98").Value = 45458476
$ws.Range("M98").Value = -45456978

$ws.Range("H122").Value = 37040788
$ws.Range("I122").Value = 45458476
$ws.Range("K122").Value = 136375524
$ws.Range("M122").Value = -136373074

$ws.Range("H128").Value = 106000
$ws.Range("J128").Value = 106000
$ws.Range("L128").Value = 106000
$ws.Range("N128").Value = -115960

$ws.Range("H132").Value = 1330.2075
$ws.Range("I132").Value = 1080.5714
$ws.Range("K132").Value = 3241.7142
$ws.Range("M132").Value = -711.7142000000003

$ws.Range("H138").Value = 3816.3257
$ws.Range("I138").Value = 881.8570999999999
$ws.Range("J138").Value = 6617.409
$ws.Range("K138").Value = 2645.5713
$ws.Range("L138").Value = 19852.227
$ws.Range("M138").Value = 2494.4287
$ws.Range("N138").Value = -30132.227

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1972.1034
$ws.Range("I88").Value = 1349.7693
$ws.Range("J88").Value = 2477.75
$ws.Range("K88").Value = 1349.7693
$ws.Range("L88").Value = 2477.75
$ws.Range("M88").Value = -943.7692999999999
$ws.Range("N88").Value = -3289.75

$ws.Range("H91").Value = 1972.1034
$ws.Range("I91").Value = 1349.7693
$ws.Range("J91").Value = 2477.75
$ws.Range("K91").Value = 1349.7693
$ws.Range("L91").Value = 2477.75
$ws.Range("M91").Value = 54.23070000000007
$ws.Range("N91").Value = -5285.75

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0

$ws.Range("H122").Value = 15961.75
$ws.Range("I122").Value = 19588.545
$ws.Range("K122").Value = 58765.63499999999
$ws.Range("M122").Value = -56315.63499999999

$ws.Range("H132").Value = 4990.4917
$ws.Range("I132").Value = 3617.1333
$ws.Range("K132").Value = 10851.3999
$ws.Range("M132").Value = -8321.3999

$ws.Range("H138").Value = 78807.664
$ws.Range("J138").Value = 78499.5
$ws.Range("L138").Value = 78499.5
$ws.Range("N138").Value = -88779.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H20").Value = 15153124
$ws.Range("I20").Value = 20835092
$ws.Range("K20").Value = 20835092
$ws.Range("M20").Value = -20834845

$ws.Range("H22").Value = 320.77777
$ws.Range("I22").Value = 298.14285
$ws.Range("K22").Value = 298.14285
$ws.Range("M22").Value = -125.14285

$ws.Range("H94").Value = 5920.222
$ws.Range("I94").Value = 2068
$ws.Range("K94").Value = 2068
$ws.Range("M94").Value = -1617

$ws.Range("H105").Value = 46564.973
$ws.Range("I105").Value = 61386.04
$ws.Range("K105").Value = 61386.04
$ws.Range("M105").Value = -59639.04

$ws.Range("H134").Value = 7070.974
$ws.Range("I134").Value = 3926.8572
$ws.Range("K134").Value = 11780.5716
$ws.Range("M134").Value = -9245.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10217.026
$ws.Range("I31").Value = 4832.5
$ws.Range("K31").Value = 4832.5
$ws.Range("M31").Value = -4537.5

$ws.Range("H34").Value = 10217.026
$ws.Range("I34").Value = 4832.5
$ws.Range("K34").Value = 4832.5
$ws.Range("M34").Value = -4630.5

$ws.Range("H134").Value = 7345.5635
$ws.Range("I134").Value = 6576.893
$ws.Range("K134").Value = 19730.679
$ws.Range("M134").Value = -17195.679

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 250.33333
$ws.Range("I26").Value = 86.666664
$ws.Range("J26").Value = 277.6111
$ws.Range("K26").Value = 259.999992
$ws.Range("L26").Value = 832.8333
$ws.Range("M26").Value = 28.00000799999998
$ws.Range("N26").Value = -1408.8333

$ws.Range("H33").Value = 261.6
$ws.Range("I33").Value = 90
$ws.Range("J33").Value = 433.2
$ws.Range("K33").Value = 540
$ws.Range("L33").Value = 2599.2
$ws.Range("M33").Value = -257
$ws.Range("N33").Value = -3165.2

$ws.Range("H40").Value = 118.5
$ws.Range("I40").Value = 70.333336
$ws.Range("J40").Value = 166.66667
$ws.Range("K40").Value = 281.333344
$ws.Range("L40").Value = 666.66668
$ws.Range("M40").Value = -212.333344
$ws.Range("N40").Value = -804.66668

$ws.Range("H80").Value = 21743564
$ws.Range("I80").Value = 4176.4707
$ws.Range("J80").Value = 83338500
$ws.Range("K80").Value = 12529.4121
$ws.Range("L80").Value = 250015500
$ws.Range("M80").Value = -11593.4121
$ws.Range("N80").Value = -250017372

$ws.Range("H83").Value = 21743564
$ws.Range("I83").Value = 4176.4707
$ws.Range("J83").Value = 83338500
$ws.Range("K83").Value = 37588.2363
$ws.Range("L83").Value = 750046500
$ws.Range("M83").Value = -32908.2363
$ws.Range("N83").Value = -750055860

$ws.Range("H86").Value = 524.8
$ws.Range("I86").Value = 533.2222
$ws.Range("J86").Value = 449
$ws.Range("K86").Value = 1599.6666
$ws.Range("L86").Value = 1347
$ws.Range("M86").Value = -413.6666
$ws.Range("N86").Value = -3719

$ws.Range("H89").Value = 524.8
$ws.Range("I89").Value = 533.2222
$ws.Range("J89").Value = 449
$ws.Range("K89").Value = 4798.999800000001
$ws.Range("L89").Value = 4041
$ws.Range("M89").Value = 1129.000199999999
$ws.Range("N89").Value = -15897

$ws.Range("H131").Value = 1732.8462
$ws.Range("I131").Value = 971.6667
$ws.Range("J131").Value = 2071.1482
$ws.Range("K131").Value = 2915.0001
$ws.Range("L131").Value = 6213.444600000001
$ws.Range("M131").Value = 2124.9999
$ws.Range("N131").Value = -16293.4446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 62273.223
$ws.Range("I122").Value = 95821
$ws.Range("K122").Value = 287463
$ws.Range("M122").Value = -285013

$ws.Range("H126").Value = 2101.889
$ws.Range("I126").Value = 1844.6
$ws.Range("J126").Value = 2423.5
$ws.Range("K126").Value = 5533.799999999999
$ws.Range("L126").Value = 7270.5
$ws.Range("M126").Value = -3063.799999999999
$ws.Range("N126").Value = -12210.5

$ws.Range("H132").Value = 4552.375
$ws.Range("I132").Value = 3662.72
$ws.Range("J132").Value = 7729.7144
$ws.Range("K132").Value = 10988.16
$ws.Range("L132").Value = 23189.1432
$ws.Range("M132").Value = -8458.16
$ws.Range("N132").Value = -28249.1432

$ws.Range("H135").Value = 86113.336
$ws.Range("J135").Value = 86113.336
$ws.Range("L135").Value = 86113.336
$ws.Range("N135").Value = -96253.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6205.2964
$ws.Range("I7").Value = 4260.6924
$ws.Range("J7").Value = 8011
$ws.Range("K7").Value = 4260.6924
$ws.Range("L7").Value = 8011
$ws.Range("M7").Value = -4148.6924
$ws.Range("N7").Value = -8235

$ws.Range("H62").Value = 40749.668
$ws.Range("J62").Value = 40749.668
$ws.Range("L62").Value = 40749.668
$ws.Range("N62").Value = -41997.668

$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15450

$ws.Range("H65").Value = 40749.668
$ws.Range("J65").Value = 40749.668
$ws.Range("L65").Value = 122249.004
$ws.Range("N65").Value = -128489.004

$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16560

$ws.Range("H80").Value = 74999
$ws.Range("J80").Value = 74999
$ws.Range("L80").Value = 74999
$ws.Range("N80").Value = -77245

$ws.Range("H83").Value = 74999
$ws.Range("J83").Value = 74999
$ws.Range("L83").Value = 224997
$ws.Range("N83").Value = -236229

$ws.Range("H126").Value = 6205.2964
$ws.Range("I126").Value = 4260.6924
$ws.Range("J126").Value = 8011
$ws.Range("K126").Value = 12782.0772
$ws.Range("L126").Value = 24033
$ws.Range("M126").Value = -10312.0772
$ws.Range("N126").Value = -28973

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0

$ws.Range("H122").Value = 214903.69
$ws.Range("I122").Value = 669036.7
$ws.Range("J122").Value = 5303.846
$ws.Range("K122").Value = 2007110.1
$ws.Range("L122").Value = 15911.538
$ws.Range("M122").Value = -2004660.1
$ws.Range("N122").Value = -20811.538

$ws.Range("H126").Value = 3600.4546
$ws.Range("I126").Value = 2622.7778
$ws.Range("K126").Value = 7868.3334
$ws.Range("M126").Value = -5398.3334

$ws.Range("H132").Value = 41692324
$ws.Range("I132").Value = 125005110
$ws.Range("J132").Value = 35930.75
$ws.Range("K132").Value = 375015330
$ws.Range("L132").Value = 107792.25
$ws.Range("M132").Value = -375012800
$ws.Range("N132").Value = -112852.25

$ws.Range("H136").Value = 41710956
$ws.Range("I136").Value = 66667308
$ws.Range("K136").Value = 200001924
$ws.Range("M136").Value = -199999374
